$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value (35 -> 48)
$ws.Range("B2").Value = 48

# Update A3 value (2 -> 1) and B3 value (13 -> 9)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 9

# Remove row 4 entirely (delete the row so dimension shrinks to A1:B3)
$ws.Rows.Item(4).Delete()
